{"js": "// The document contains three \"<id>p009r_N</id>\" markers (N = 1, 2, 3),\n// each one originally split across three separate runs:\n//   1) \"<id>\"        (Courier New, color 7f6000, sz 18)\n//   2) \"p009r_N\"     (color 000000, default font)\n//   3) \"</id>\"       (Courier New, color 7f6000, sz 18)\n//\n// The edit merges each trio of runs into a single run whose text is the\n// concatenation \"<id>p009r_N</id>\" and whose formatting is that of the\n// first (\"<id>\") run. Re-inserting the found range's text with\n// InsertLocation.Replace collapses the matched run(s) into one run that\n// carries forward the formatting of the range's leading run, which is\n// exactly the desired result.\n\nconst body = context.document.body;\n\nfor (const n of [1, 2, 3]) {\n  const needle = `<id>p009r_${n}</id>`;\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    // Nothing to do for this id if it can't be found (already merged, etc).\n    continue;\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    const found = results.items[i];\n    found.insertText(needle, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three \"<id>p009r_N</id>\" markers (N = 1, 2, 3),\n# each one originally split across three separate runs:\n#   1) \"<id>\"        (Courier New, color 7f6000, sz 18)\n#   2) \"p009r_N\"     (color 000000, default font)\n#   3) \"</id>\"       (Courier New, color 7f6000, sz 18)\n#\n# The edit merges each trio of runs into a single run whose text is the\n# concatenation \"<id>p009r_N</id>\" and whose formatting is that of the\n# first (\"<id>\") run. Running Find/Replace (wdReplaceAll) over that exact\n# literal text rewrites the matched span as a single run carrying the\n# formatting of the first run in the match, which is exactly what we want.\n\n$d = $word.ActiveDocument\n\nforeach ($n in 1..3) {\n    $needle = \"<id>p009r_$n</id>\"\n\n    $rng = $d.Content\n    # 1=FindText, 2=MatchCase, 3=MatchWholeWord, 4=MatchWildcards,\n    # 5=MatchSoundsLike, 6=MatchAllWordForms, 7=Forward, 8=Wrap(wdFindContinue=1),\n    # 9=Format, 10=ReplaceWith, 11=Replace(wdReplaceAll=2)\n    $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null\n}\n"}
